$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 165445
$ws.Range("C4").Value = 156391
$ws.Range("C7").Value = 5.47
$ws.Range("C8").Value = 65.08
